# Commit: "Generate Report for Archive"
#
# The localization record for "a7e0052c-a61f-42c9-8f6d-436456ce933a" moved
# from the last data row (row 9) up to row 5 on every sheet (Overview,
# zh-cn, de-de), pushing the records that used to occupy rows 5-8
# (bf9987fa..., ed5c1042..., 1fbda1af..., 52679a15...) down by one row.
#
# We reproduce that by rewriting the cell values of rows 5-9 in place
# (the header row 1 and the first three data rows 2-4 are untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn), C (de-de),
# D (Latest Handoff Date)
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ovRows = @(
    @("a7e0052c-a61f-42c9-8f6d-436456ce933a.md", "Ready for handoff", "Ready for handoff", "2016-36-13 00:36:37"),
    @("bf9987fa-933a-4d56-a631-b55c9c97b021.md", "In Translation",    "In Translation",    "2016-34-13 00:34:21"),
    @("ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "In Translation",    "In Translation",    "2016-33-13 00:33:49"),
    @("1fbda1af-7bda-4ec9-8990-163d2f0d8608.md", "Ready for handoff", "Ready for handoff", "2016-32-13 00:32:46"),
    @("52679a15-f857-4c26-9a01-c476e58b1a39.md", "Ready for handoff", "Ready for handoff", "2016-34-13 00:34:37")
)

for ($i = 0; $i -lt $ovRows.Length; $i++) {
    $r = 5 + $i
    $row = $ovRows[$i]
    $ov.Cells.Item($r, 1).Value = $row[0]
    $ov.Cells.Item($r, 2).Value = $row[1]
    $ov.Cells.Item($r, 3).Value = $row[2]
    $ov.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------
# Language detail sheets: columns A (Source File Name),
# B (File Extension, unchanged ".md"), C (Status),
# D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------
function Set-LangRows($ws, $rows) {
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = 5 + $i
        $row = $rows[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 3).Value = $row[1]
        $ws.Cells.Item($r, 4).Value = $row[2]
        $ws.Cells.Item($r, 5).Value = $row[3]
    }
}

$zh = $wb.Worksheets.Item("zh-cn")
$zhRows = @(
    @("a7e0052c-a61f-42c9-8f6d-436456ce933a.md", "Ready for handoff", "a7e0052c-a61f-42c9-8f6d-436456ce933a.0855e895c9696c904afa5cffcf2284f4df0087da.zh-cn.xlf", "2016-03-13 00:36:34"),
    @("bf9987fa-933a-4d56-a631-b55c9c97b021.md", "In Translation",    "bf9987fa-933a-4d56-a631-b55c9c97b021.ddfe39798ce3afc4d6c2a625a951d48407f3e769.zh-cn.xlf", "2016-03-13 00:34:17"),
    @("ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "In Translation",    "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.zh-cn.xlf", "2016-03-13 00:33:45"),
    @("1fbda1af-7bda-4ec9-8990-163d2f0d8608.md", "Ready for handoff", "1fbda1af-7bda-4ec9-8990-163d2f0d8608.683fa9eb4c7f22a59dff9a742a354e862bcc2f6d.zh-cn.xlf", "2016-03-13 00:32:42"),
    @("52679a15-f857-4c26-9a01-c476e58b1a39.md", "Ready for handoff", "52679a15-f857-4c26-9a01-c476e58b1a39.62d1bde33acbc6125c05a01f51b09012a18dd835.zh-cn.xlf", "2016-03-13 00:34:34")
)
Set-LangRows $zh $zhRows

$de = $wb.Worksheets.Item("de-de")
$deRows = @(
    @("a7e0052c-a61f-42c9-8f6d-436456ce933a.md", "Ready for handoff", "a7e0052c-a61f-42c9-8f6d-436456ce933a.0855e895c9696c904afa5cffcf2284f4df0087da.de-de.xlf", "2016-03-13 00:36:37"),
    @("bf9987fa-933a-4d56-a631-b55c9c97b021.md", "In Translation",    "bf9987fa-933a-4d56-a631-b55c9c97b021.ddfe39798ce3afc4d6c2a625a951d48407f3e769.de-de.xlf", "2016-03-13 00:34:21"),
    @("ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "In Translation",    "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.de-de.xlf", "2016-03-13 00:33:49"),
    @("1fbda1af-7bda-4ec9-8990-163d2f0d8608.md", "Ready for handoff", "1fbda1af-7bda-4ec9-8990-163d2f0d8608.683fa9eb4c7f22a59dff9a742a354e862bcc2f6d.de-de.xlf", "2016-03-13 00:32:46"),
    @("52679a15-f857-4c26-9a01-c476e58b1a39.md", "Ready for handoff", "52679a15-f857-4c26-9a01-c476e58b1a39.62d1bde33acbc6125c05a01f51b09012a18dd835.de-de.xlf", "2016-03-13 00:34:37")
)
Set-LangRows $de $deRows
